$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 44 (old rows 44-46 -> 45-47, SUM()/formula refs auto-adjust).
$ws.Rows("44:44").Insert()

# Mark rows 37-43 column A with the "yellow" highlight fill (reuses existing style s="9").
$ws.Range("A37:A44").Interior.Color = 65535

# Move the part-number text that lived in column D (rows 38-43) into column C instead.
for ($r = 38; $r -le 43; $r++) {
    $partNo = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $partNo
    $ws.Cells.Item($r, 4).ClearContents()
}

# Fill in the new row 44 with the extra purchased part: APT15DQ60KG.
$ws.Range("C44").Value = "APT15DQ60KG"
$ws.Range("F44").Value = "APT15DQ60KG"
$ws.Range("F44").Interior.Color = 65535
$ws.Range("G44").Value = 4.49
$ws.Range("N44").Value = "https://www.digikey.se/products/sv?keywords=APT15DQ60KG"

# Hyperlink the purchase link in N44, then restore its (non-default) text style.
$ws.Hyperlinks.Add($ws.Cells.Item(44, 14), "https://www.digikey.se/products/sv?keywords=APT15DQ60KG")
$ws.Range("N44").Style = "Hyperlänk"

# Update selection/viewport to match where the author left off editing.
$ws.Range("B39").Select()

Write-Output "done"
